$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the existing header cell (G1) onto the new
# header cell (H1) so it shares the same bold/border/centered style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column, one per data row (2-13)
$saveValues = @(1, 1, 0, 0, 0, 1, 0, 0, 1, 0, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
